$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLink = "https://www.360dx.com/business-news/top-five-articles-360dx-last-week-fda-aims-reclassify-cdx-tests-tempus-buying-oneome"
$newKeywords = "CDx"
$newTitle = "Top Five Articles on 360Dx Last Week: FDA Aims to Reclassify CDx Tests; Tempus Buying OneOme Assets; More"

# New row goes right after the current last data row (row 67) -> row 68
$row = 68

$ws.Range("B$row").Value = $newKeywords
$ws.Range("C$row").Value = $newTitle

# Create the hyperlink (this also sets the cell text/value + registers the
# external relationship + <hyperlinks> entry).
$ws.Hyperlinks.Add($ws.Range("A$row"), $newLink)

# Match the same visual style used by the other link cells in column A
# (Hyperlinks.Add leaves behind a slightly different auto-generated style).
$ws.Range("A$row").Style = $ws.Range("A67").Style
